# Updated symbol list on Wed Dec 28 23:43:09 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numeric-looking values as plain text
# (inlineStr in the source file, no numeric cell format). Assigning a
# numeric-looking string straight to .Value lets Excel auto-convert it to
# a real number, which would change the cell's stored type. Prefixing
# with an apostrophe forces Excel to keep it as text (quote-prefixed)
# and ClearFormats() afterwards drops the quote-prefix cell style that
# gets stamped on, restoring the cell to its original (default/general)
# formatting while keeping the text value intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.ClearFormats()
}

# Price (column D) updates
Set-TextValue "D3"  "23.85"
Set-TextValue "D4"  "5.234"
Set-TextValue "D5"  "0.05764"
Set-TextValue "D6"  "6.409"
Set-TextValue "D7"  "3.241"
Set-TextValue "D8"  "0.8111"
Set-TextValue "D9"  "0.8834"
Set-TextValue "D10" "0.1371"
Set-TextValue "D11" "0.06971"
Set-TextValue "D12" "0.03172"
Set-TextValue "D13" "0.03044"
Set-TextValue "D14" "0.09322"
Set-TextValue "D15" "3.810"
Set-TextValue "D16" "0.001521"
Set-TextValue "D17" "0.04697"
Set-TextValue "D18" "0.0006033"
Set-TextValue "D19" "0.006171"
Set-TextValue "D21" "0.004070"
Set-TextValue "D22" "0.00008693"
Set-TextValue "D23" "3.546"
Set-TextValue "D24" "2.146"
Set-TextValue "D40" "0.03729"
Set-TextValue "D41" "0.006245"
Set-TextValue "D43" "0.002403"
Set-TextValue "D44" "0.007953"
Set-TextValue "D45" "0.00005294"
Set-TextValue "D48" "0.002412"
Set-TextValue "D49" "0.00002098"

# Volume(1h) (column E) label updates — "Bestin24h" / "Worstin24h" badges
# being added/removed/relabeled alongside the price refresh
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
